$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B8").NumberFormat = "@"

$ws.Range("A8").Value = "2025-11-11 03:57:57"
$ws.Range("B8").Value = "2023-03-20"
$ws.Range("C8").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("D8").Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/HZL20032023.pdf"

$ws.Range("B8").ClearFormats()
